$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New TPM-derived values for the Cd22-Ptprc LR-pair sheet.
# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E..T numeric metrics. A new target cluster "MuSCs" was added, so each of the
# three sending clusters now has 6 target-cluster rows instead of 5 (rows 2-19).
$data = @(
    @("Inflammatory-Mac", "Cd22", "Ptprc", "ECs", 3, 1, "1.273363333333333", "3.82009", "0.1769428433887536", "0.1769428433887536", 2, 1, "0.384696", "0.7693920000000001", "0.0003953339325349752", "0.000263641226651746", "0.4898577808800001", "2.939146685280001", "6.995151011079618e-05", "4.664942827825878e-05"),
    @("Inflammatory-Mac", "Cd22", "Ptprc", "FAPs", 3, 1, "1.273363333333333", "3.82009", "0.1769428433887536", "0.1769428433887536", 1, "0.3333333333333333", "0.005387666666666666", "0.016163", "5.536650889328026e-06", "5.538442232791829e-06", "0.006860457185555556", "0.06174411467000001", "9.796707512085721e-07", "9.799877166145432e-07"),
    @("Inflammatory-Mac", "Cd22", "Ptprc", "Inflammatory-Mac", 3, 1, "1.273363333333333", "3.82009", "0.1769428433887536", "0.1769428433887536", 3, 1, "272.782435", "818.347305", "0.2803256409087078", "0.2804163382480588", "347.3511507063834", "3126.16035635745", "0.04960161597716145", "0.04961766422227402"),
    @("Inflammatory-Mac", "Cd22", "Ptprc", "MuSCs", 3, 1, "1.273363333333333", "3.82009", "0.1769428433887536", "0.1769428433887536", 1, "0.5", "0.5595085", "1.119017", "0.0005749804926272826", "0.0003834443489458648", "0.7124576085883333", "4.27474565153", "0.0001017386832585376", "6.784773338383073e-05"),
    @("Inflammatory-Mac", "Cd22", "Ptprc", "Neutrophils", 3, 1, "1.273363333333333", "3.82009", "0.1769428433887536", "0.1769428433887536", 3, 1, "585.2320353333333", "1755.696106", "0.6014153564730961", "0.6016099400741544", "745.2130152855045", "6706.917137569541", "0.1064161432320104", "0.1064505734076585"),
    @("Neutrophils", "Cd22", "Ptprc", "Resolving-Mac", 3, 1, "1.273363333333333", "3.82009", "0.1769428433887536", "0.1769428433887536", 3, 1, "114.1272113333333", "342.381634", "0.1172831515421446", "0.1173210976599564", "145.3254062474511", "1307.92865622706", "0.02075241431546114", "0.02075912860944232"),
    @("Neutrophils", "Cd22", "Ptprc", "ECs", 3, 1, "0.815617", "2.446851", "0.1133357521128076", "0.1133357521128076", 2, 1, "0.384696", "0.7693920000000001", "0.0003953339325349752", "0.000263641226651746", "0.313764597432", "1.882587584592", "4.480546857956534e-05", "2.98799767105188e-05"),
    @("Neutrophils", "Cd22", "Ptprc", "FAPs", 3, 1, "0.815617", "2.446851", "0.1133357521128076", "0.1133357521128076", 1, "0.3333333333333333", "0.005387666666666666", "0.016163", "5.536650889328026e-06", "5.538442232791829e-06", "0.004394272523666667", "0.039548452713", "6.275004927280367e-07", "6.277035159867991e-07"),
    @("Neutrophils", "Cd22", "Ptprc", "Inflammatory-Mac", 3, 1, "0.815617", "2.446851", "0.1133357521128076", "0.1133357521128076", 3, 1, "272.782435", "818.347305", "0.2803256409087078", "0.2804163382480588", "222.485991287395", "2002.373921586555", "0.03177091734889321", "0.0317811966000632"),
    @("Neutrophils", "Cd22", "Ptprc", "MuSCs", 3, 1, "0.815617", "2.446851", "0.1133357521128076", "0.1133357521128076", 1, "0.5", "0.5595085", "1.119017", "0.0005749804926272826", "0.0003834443489458648", "0.4563446442445", "2.738067865467", "6.516584658210568e-05", "4.345795368118542e-05"),
    @("Neutrophils", "Cd22", "Ptprc", "Neutrophils", 3, 1, "0.815617", "2.446851", "0.1133357521128076", "0.1133357521128076", 3, 1, "585.2320353333333", "1755.696106", "0.6014153564730961", "0.6016099400741544", "477.3251969624674", "4295.926772662206", "0.06816186175807061", "0.06818391503684537"),
    @("Neutrophils", "Cd22", "Ptprc", "Resolving-Mac", 3, 1, "0.815617", "2.446851", "0.1133357521128076", "0.1133357521128076", 3, 1, "114.1272113333333", "342.381634", "0.1172831515421446", "0.1173210976599564", "93.08409372605934", "837.7568435345341", "0.01329237419018934", "0.0132966748419913"),
    @("Resolving-Mac", "Cd22", "Ptprc", "ECs", 3, 1, "5.107486666666667", "15.32246", "0.7097214044984388", "0.7097214044984388", 2, 1, "0.384696", "0.7693920000000001", "0.0003953339325349752", "0.000263641226651746", "1.96482969072", "11.78897814432", "0.0002805769538446136", "0.0001871118216629684"),
    @("Resolving-Mac", "Cd22", "Ptprc", "FAPs", 3, 1, "5.107486666666667", "15.32246", "0.7097214044984388", "0.7097214044984388", 1, "0.3333333333333333", "0.005387666666666666", "0.016163", "5.536650889328026e-06", "5.538442232791829e-06", "0.02751743566444445", "0.24765692098", "3.929479645391417e-06", "3.930751000190486e-06"),
    @("Resolving-Mac", "Cd22", "Ptprc", "Inflammatory-Mac", 3, 1, "5.107486666666667", "15.32246", "0.7097214044984388", "0.7097214044984388", 3, 1, "272.782435", "818.347305", "0.2803256409087078", "0.2804163382480588", "1393.232649663367", "12539.0938469703", "0.1989531075826531", "0.1990174774257216"),
    @("Resolving-Mac", "Cd22", "Ptprc", "MuSCs", 3, 1, "5.107486666666667", "15.32246", "0.7097214044984388", "0.7097214044984388", 1, "0.5", "0.5595085", "1.119017", "0.0005749804926272826", "0.0003834443489458648", "2.857682203636667", "17.14609322182", "0.0004080759627866392", "0.0002721386618808486"),
    @("Resolving-Mac", "Cd22", "Ptprc", "Neutrophils", 3, 1, "5.107486666666667", "15.32246", "0.7097214044984388", "0.7097214044984388", 3, 1, "585.2320353333333", "1755.696106", "0.6014153564730961", "0.6016099400741544", "2989.064817371196", "26901.58335634076", "0.426837351483015", "0.4269754516296504"),
    @("Resolving-Mac", "Cd22", "Ptprc", "Resolving-Mac", 3, 1, "5.107486666666667", "15.32246", "0.7097214044984388", "0.7097214044984388", 3, 1, "114.1272113333333", "342.381634", "0.1172831515421446", "0.1173210976599564", "582.903210188849", "5246.128891699641", "0.0832383630364941", "0.08326529420852272")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    for ($c = 5; $c -le 20; $c++) {
        $ws.Cells.Item($r, $c).Value = [double]$row[$c-1]
    }
    $r = $r + 1
}
